$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.98"
$ws.Range("E2").Value = "'-0.89%"

$ws.Range("D3").Value = "'37.42"
$ws.Range("E3").Value = "'-0.44%"

$ws.Range("D4").Value = "'5.124"
$ws.Range("E4").Value = "'1.03%"

$ws.Range("D5").Value = "'0.07797"

$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.435"
$ws.Range("E6").Value = "'1.88%"

$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.251"
$ws.Range("E7").Value = "'0.33%"

$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'1.889"
$ws.Range("E8").Value = "'0.43%"

$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.988"
$ws.Range("E9").Value = "'2.61%"

$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9264"
$ws.Range("E10").Value = "'0.56%"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1082"
$ws.Range("E11").Value = "'-10.03%"

$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1919"
$ws.Range("E12").Value = "'-0.21%"

$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.08953"
$ws.Range("E13").Value = "'-4.08%"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03328"
$ws.Range("E14").Value = "'-3.08%"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09584"
$ws.Range("E15").Value = "'-0.94%"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001377"
$ws.Range("E16").Value = "'0.47%"

$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005757"
$ws.Range("E17").Value = "'-1.16%"

$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.552"
$ws.Range("E18").Value = "'-0.03%"

$ws.Range("D19").Value = "'0.3476"
$ws.Range("E19").Value = "'2.18%"

$ws.Range("E20").Value = "'18.47%"

$ws.Range("E21").Value = "'-1.69%"

$ws.Range("D22").Value = "'0.2585"
$ws.Range("E22").Value = "'-0.22%"

$ws.Range("E23").Value = "'0.71%"

$ws.Range("E24").Value = "'-0.91%"

$ws.Range("D25").Value = "'0.004252"
$ws.Range("E25").Value = "'-0.04%"

$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'0.14%"

$ws.Range("D39").Value = "'0.02173"
$ws.Range("E39").Value = "'3.91%"

$ws.Range("D40").Value = "'0.05026"
$ws.Range("E40").Value = "'0.97%"

$ws.Range("D41").Value = "'0.007445"
$ws.Range("E41").Value = "'-2.55%"

$ws.Range("E42").Value = "'0.25%"

$ws.Range("D43").Value = "'0.008662"
$ws.Range("E43").Value = "'-11.90%"

$ws.Range("D44").Value = "'0.002114"
$ws.Range("E44").Value = "'2.57%"

$ws.Range("E45").Value = "'-9.51%"

$ws.Range("D46").Value = "'0.00006533"
$ws.Range("E46").Value = "'-1.74%"

$ws.Range("E47").Value = "'0.14%"

$ws.Range("D48").Value = "'0.002873"
$ws.Range("E48").Value = "'-2.14%"

$ws.Range("E49").Value = "'-16.55%"

$ws.Range("E50").Value = "'0.14%"

$ws.Range("E51").Value = "'0.14%"
